# "Raw Matrix" holds the hand-entered (un-normalized) Markov transition
# weights; "Normalized Matrix" derives each row by dividing by the row sum,
# so editing one raw weight automatically ripples through the normalized
# sheet's formulas on recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Raw Matrix")
$ws.Activate()

# Raw weight for the "Coffee Grinder" -> "Coffee Grinder" transition in row 14
# (First Position) was lowered from 0.8 to 0.3.
$ws.Range("N14").Value = 0.3

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("N15").Select()
